# Auto-generated Excel COM-interop edit script
# Applies cached-value corrections to Leve price/profit columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the commit diff.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 198.25
$ws.Range("I5").Value = 223.71428
$ws.Range("K5").Value = 223.71428
$ws.Range("M5").Value = -108.71428
$ws.Range("H6").Value = 448.66666
$ws.Range("I6").Value = 448.66666
$ws.Range("K6").Value = 1345.99998
$ws.Range("M6").Value = -1233.99998
$ws.Range("H8").Value = 154.58333
$ws.Range("I8").Value = 154.58333
$ws.Range("K8").Value = 463.74999
$ws.Range("M8").Value = -324.74999
$ws.Range("H17").Value = 1851.8684
$ws.Range("J17").Value = 1851.8684
$ws.Range("L17").Value = 5555.6052
$ws.Range("N17").Value = -5891.6052
$ws.Range("H20").Value = 3997.5
$ws.Range("I20").Value = 3997.5
$ws.Range("K20").Value = 3997.5
$ws.Range("M20").Value = -3767.5
$ws.Range("H33").Value = 374.35715
$ws.Range("I33").Value = 387.92307
$ws.Range("K33").Value = 387.92307
$ws.Range("M33").Value = -158.92307
$ws.Range("H35").Value = 3997.5
$ws.Range("I35").Value = 3997.5
$ws.Range("K35").Value = 3997.5
$ws.Range("M35").Value = -3618.5
$ws.Range("H44").Value = 23331.666
$ws.Range("J44").Value = 23331.666
$ws.Range("L44").Value = 23331.666
$ws.Range("N44").Value = -24255.666
$ws.Range("H74").Value = 4993
$ws.Range("I74").Value = 4993
$ws.Range("K74").Value = 4993
$ws.Range("M74").Value = -4057
$ws.Range("H77").Value = 4993
$ws.Range("I77").Value = 4993
$ws.Range("K77").Value = 24965
$ws.Range("M77").Value = -20285
$ws.Range("H86").Value = 5225
$ws.Range("I86").Value = 5000
$ws.Range("J86").Value = 5450
$ws.Range("K86").Value = 5000
$ws.Range("L86").Value = 5450
$ws.Range("M86").Value = -3877
$ws.Range("N86").Value = -7696
$ws.Range("H89").Value = 5225
$ws.Range("I89").Value = 5000
$ws.Range("J89").Value = 5450
$ws.Range("K89").Value = 25000
$ws.Range("L89").Value = 27250
$ws.Range("M89").Value = -19384
$ws.Range("N89").Value = -38482
$ws.Range("H107").Value = 1657.75
$ws.Range("I107").Value = 1663.5
$ws.Range("J107").Value = 1640.5
$ws.Range("K107").Value = 1663.5
$ws.Range("L107").Value = 1640.5
$ws.Range("M107").Value = 256.5
$ws.Range("N107").Value = -5480.5

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 35285.43
$ws.Range("I32").Value = 35285.43
$ws.Range("K32").Value = 35285.43
$ws.Range("M32").Value = -34998.43
$ws.Range("H80").Value = 57696.5
$ws.Range("I80").Value = 34000
$ws.Range("J80").Value = 67852.14
$ws.Range("K80").Value = 34000
$ws.Range("L80").Value = 67852.14
$ws.Range("M80").Value = -33002
$ws.Range("N80").Value = -69848.14
$ws.Range("H83").Value = 57696.5
$ws.Range("I83").Value = 34000
$ws.Range("J83").Value = 67852.14
$ws.Range("K83").Value = 102000
$ws.Range("L83").Value = 203556.42
$ws.Range("M83").Value = -97008
$ws.Range("N83").Value = -213540.42
$ws.Range("H88").Value = 1390.8462
$ws.Range("I88").Value = 965.6667
$ws.Range("J88").Value = 1755.2858
$ws.Range("K88").Value = 965.6667
$ws.Range("L88").Value = 1755.2858
$ws.Range("M88").Value = -559.6667
$ws.Range("N88").Value = -2567.2858
$ws.Range("H91").Value = 1390.8462
$ws.Range("I91").Value = 965.6667
$ws.Range("J91").Value = 1755.2858
$ws.Range("K91").Value = 965.6667
$ws.Range("L91").Value = 1755.2858
$ws.Range("M91").Value = 438.3333
$ws.Range("N91").Value = -4563.2858
$ws.Range("H122").Value = 2687.5386
$ws.Range("I122").Value = 2687.5386
$ws.Range("K122").Value = 8062.6158
$ws.Range("M122").Value = -5612.6158
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2553.0527
$ws.Range("I86").Value = 2639.389
$ws.Range("K86").Value = 2639.389
$ws.Range("M86").Value = -1516.389
$ws.Range("H89").Value = 2553.0527
$ws.Range("I89").Value = 2639.389
$ws.Range("K89").Value = 13196.945
$ws.Range("M89").Value = -7580.945
$ws.Range("H105").Value = 3953.3157
$ws.Range("I105").Value = 2357.3333
$ws.Range("K105").Value = 2357.3333
$ws.Range("M105").Value = -610.3332999999998

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 599.6667
$ws.Range("I2").Value = 599.6667
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 599.6667
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -486.6667
$ws.Range("N2").ClearContents()
$ws.Range("H22").Value = 10623.385
$ws.Range("I22").Value = 1933.3334
$ws.Range("J22").Value = 18072
$ws.Range("K22").Value = 1933.3334
$ws.Range("L22").Value = 18072
$ws.Range("M22").Value = -1583.3334
$ws.Range("N22").Value = -18772
$ws.Range("H62").Value = 2777.4
$ws.Range("I62").Value = 2499.5
$ws.Range("K62").Value = 2499.5
$ws.Range("M62").Value = -1875.5
$ws.Range("H65").Value = 2777.4
$ws.Range("I65").Value = 2499.5
$ws.Range("K65").Value = 12497.5
$ws.Range("M65").Value = -9377.5
$ws.Range("H134").Value = 3293.5833
$ws.Range("I134").Value = 3343
$ws.Range("K134").Value = 10029
$ws.Range("M134").Value = -7494

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 17058.111
$ws.Range("J7").Value = 74
$ws.Range("L7").Value = 222
$ws.Range("N7").Value = -446
$ws.Range("H13").Value = 650
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H22").Value = 100
$ws.Range("J22").Value = 100
$ws.Range("L22").Value = 300
$ws.Range("N22").Value = -638
$ws.Range("H23").Value = 495.16666
$ws.Range("J23").Value = 693
$ws.Range("L23").Value = 2079
$ws.Range("N23").Value = -2549
$ws.Range("H27").Value = 100
$ws.Range("J27").Value = 100
$ws.Range("L27").Value = 300
$ws.Range("N27").Value = -504
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2000
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 2000
$ws.Range("L80").Value = 2000
$ws.Range("M80").Value = -1002
$ws.Range("N80").Value = -3996
$ws.Range("H83").Value = 2000
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 10000
$ws.Range("L83").Value = 10000
$ws.Range("M83").Value = -5008
$ws.Range("N83").Value = -19984

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4720.125
$ws.Range("J7").Value = 3005
$ws.Range("L7").Value = 3005
$ws.Range("N7").Value = -3229
$ws.Range("H12").Value = 2016.25
$ws.Range("J12").Value = 2016.25
$ws.Range("L12").Value = 2016.25
$ws.Range("N12").Value = -2356.25
$ws.Range("H22").Value = 2702.5715
$ws.Range("J22").Value = 2780
$ws.Range("L22").Value = 2780
$ws.Range("N22").Value = -3370
$ws.Range("H27").Value = 2702.5715
$ws.Range("J27").Value = 2780
$ws.Range("L27").Value = 2780
$ws.Range("N27").Value = -2994
$ws.Range("H46").Value = 3307.1667
$ws.Range("I46").Value = 2197.25
$ws.Range("J46").Value = 5527
$ws.Range("K46").Value = 2197.25
$ws.Range("L46").Value = 5527
$ws.Range("M46").Value = -2009.25
$ws.Range("N46").Value = -5903
$ws.Range("H55").Value = 1360.6428
$ws.Range("I55").Value = 1175
$ws.Range("J55").Value = 1434.9
$ws.Range("K55").Value = 1175
$ws.Range("L55").Value = 1434.9
$ws.Range("M55").Value = -1002
$ws.Range("N55").Value = -1780.9
$ws.Range("H68").Value = 4435.6
$ws.Range("I68").Value = 4435.6
$ws.Range("K68").Value = 4435.6
$ws.Range("M68").Value = -3686.6
$ws.Range("H71").Value = 4435.6
$ws.Range("I71").Value = 4435.6
$ws.Range("K71").Value = 22178
$ws.Range("M71").Value = -18434
$ws.Range("H100").Value = 2999
$ws.Range("I100").Value = 2999
$ws.Range("K100").Value = 2999
$ws.Range("M100").Value = -2458
$ws.Range("H126").Value = 4720.125
$ws.Range("J126").Value = 3005
$ws.Range("L126").Value = 9015
$ws.Range("N126").Value = -13955
$ws.Range("H132").Value = 3837.7778
$ws.Range("I132").Value = 3863.2856
$ws.Range("K132").Value = 11589.8568
$ws.Range("M132").Value = -9059.856800000001

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1383
$ws.Range("I81").Value = 1383
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2766
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -1705
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 1383
$ws.Range("I84").Value = 1383
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 13830
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -8526
$ws.Range("N84").ClearContents()
$ws.Range("H122").Value = 2721.4285
$ws.Range("I122").Value = 3009.0908
$ws.Range("K122").Value = 9027.2724
$ws.Range("M122").Value = -6577.2724
$ws.Range("H126").Value = 2726.7856
$ws.Range("J126").Value = 3249.5
$ws.Range("L126").Value = 9748.5
$ws.Range("N126").Value = -14688.5

Write-Output "Applied Phantom_Profits corrections to 8 sheets."